$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header row (row 1) renames
# ---------------------------------------------------------------------------
$ws.Range("C1").Value = "servico_amostra"
$ws.Range("K1").Value = "id_veiculo_apurado"
$ws.Range("L1").Value = "servico_apurado"
$ws.Range("M1").Value = "sentido_apurado"
$ws.Range("N1").Value = "datetime_partida_apurado"
$ws.Range("O1").Value = "datetime_chegada_apurado"

# ---------------------------------------------------------------------------
# Data rows: the "status" message in column J changes to a new message, and
# the previously swapped servico/id_veiculo values in columns K and L are
# corrected back (K: 844 -> 47689, L: 47689 -> 844). K/L must remain text.
# ---------------------------------------------------------------------------
$newStatus = "Viagem circular identificada e já paga"

$rowsToUpdate = @(
    18,94,112,138,139,140,141,142,143,144,145,146,147,148,149,150,
    151,152,153,154,155,156,157,158,159,160,161,162,163,164,165,166,
    167,168,169,170,171,172,173,174,175,176,177,178,179,180,181,182,
    183,184,185,186,188,189,190,191,192,193,194,195,196,197,198,199,
    200,201,202,203,204,205,206,207,208,209,210,211,212,213,214,215,
    216,217,218,219,220,221,222,223,224,225,226,227,228,229,230,231,
    232,233,234,235,236,238,239,240,241,242,243,244,245,246,247,248,
    249,250,251,252,253,254,255,256,257,258,259,260,261,262,263,264,
    265,266,267,268,269,270,271,272,273,274,275,276,277,278,279,280
)

foreach ($r in $rowsToUpdate) {
    $ws.Range("J$r").Value = $newStatus

    $kCell = $ws.Range("K$r")
    $kCell.NumberFormat = "@"
    $kCell.Value = "47689"

    $lCell = $ws.Range("L$r")
    $lCell.NumberFormat = "@"
    $lCell.Value = "844"
}

# ---------------------------------------------------------------------------
# Rows 187 and 281: trim the trailing clause off the "sem sinal" message.
# ---------------------------------------------------------------------------
$trimmedMsg = "Viagem circular inválida - sem sinal inicial/final dentro do raio de 500m"
$ws.Range("J187").Value = $trimmedMsg
$ws.Range("J281").Value = $trimmedMsg
